# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 on Sheet1 with newly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.339.79"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.932.87"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'" + "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'" + "250.41"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'" + "0.7370"
$ws.Range("D7").Value = "'" + "1.001"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'" + "0.3218"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").Value = "'" + "27.88"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").Value = "'" + "0.07102"
$ws.Range("E10").Value = "  -4.55%  "
$ws.Range("D11").Value = "'" + "0.7871"
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("D12").Value = "'" + "0.08030"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "1.932.93"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "'" + "5.391"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "'" + "94.75"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "30.340.47"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'" + "254.74"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'" + "0.000008055"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "'" + "5.746"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "2.186.57"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'" + "1.001"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'" + "1.001"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'" + "6.830"
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("D25").Value = "'" + "9.563"
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("D26").Value = "'" + "164.12"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "'" + "19.10"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'" + "2.296"
$ws.Range("E28").Value = "  -5.49%  "
$ws.Range("D29").Value = "'" + "0.1324"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'" + "1.360"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("D32").Value = "'" + "4.424"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").Value = "'" + "4.155"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").Value = "'" + "0.05119"
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("D35").Value = "'" + "1.289"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("D36").Value = "'" + "0.7471"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").Value = "'" + "2.770"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "'" + "2.803"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").Value = "'" + "78.27"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("D41").Value = "'" + "6.411"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("D42").Value = "'" + "0.4508"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "'" + "1.989"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("D44").Value = "'" + "0.8449"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'" + "1.001"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'" + "101.38"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").Value = "'" + "7.557"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "'" + "9.747"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "'" + "982.70"
$ws.Range("E49").Value = "  +11.19%  "
$ws.Range("D50").Value = "'" + "37.06"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "'" + "0.4186"
$ws.Range("E51").Value = "  -0.74%  "
